$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:L1): bold font, thin box border, centered/top aligned ---
$header = $ws.Range("A1:L1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

# --- New data row 2 ---
$ws.Range("A2").Value = 8192984600
$ws.Range("B2").Value = "user_8192984600"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = "2026-02-13T20:43:55.614256+00:00"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = $false
$ws.Range("I2").Value = $false
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "Added during extraction"

# L2 looks like a date ("2026-02-13") - force it to stay plain text instead
# of being auto-converted into a date serial by prefixing with an apostrophe,
# then resetting the cell style back to Normal so no extra number format
# sticks around on the cell.
$ws.Range("L2").Value = "'2026-02-13"
$ws.Range("L2").Style = "Normal"

Write-Output "done"
